$d = $word.ActiveDocument

# Locate the paragraph that holds the inline picture (the "GATO1.bmp" drawing) -
# the new paragraph about caperucita must be inserted right after it.
$shape = $d.InlineShapes.Item(1)
$picPara = $shape.Range.Paragraphs.Item(1)
$picStart = $picPara.Range.Start

# Figure out this paragraph's 1-based index in the document's Paragraphs
# collection (more reliable than chaining .Next() across a mutation).
$picIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $picStart) {
        $picIndex = $i
    }
}

# Insert a fresh empty paragraph right after the picture's paragraph.
$picPara.Range.InsertParagraphAfter() | Out-Null

# That freshly inserted paragraph is now the very next paragraph.
$newPara = $d.Paragraphs.Item($picIndex + 1)
$insertionRange = $newPara.Range
$insertionRange.Collapse(1)

# Build the new paragraph's content exactly, including the proofing-error
# markers Word leaves around the words it has flagged, by inserting raw OOXML.
$paraXml = '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>caperucita</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> siguió su camino contenta  llevando su canasto de manzanas no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sabia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que el malvado lobo ya se había comido a la abuelita y se había puesto el traje de la abuela para que cuando llegara caperucita no lo reconociera </w:t></w:r></w:p>'

$insertionRange.InsertXML($paraXml) | Out-Null
